$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in the existing "Жо Бидон Амэрикович" row's e-mail (row 4) ---
# The old hyperlink (3rd one, pointing at C4) has to be replaced so its
# mailto target is corrected too - not just the visible text.
$idx = 0
foreach ($h in $ws.Hyperlinks) {
    $idx = $idx + 1
    if ($idx -eq 3) {
        $h.Delete()
    }
}
$ws.Range("C4").Value = "bidonchick_joe@mail.ru"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:bidonchick_joe@mail.ru")
$ws.Range("C4").Style = $ws.Range("C2").Style

# --- Append two new rows of client data ---
$ws.Range("A6").Value = "Пончик"
$ws.Range("A7").Value = "Не пончик"
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:ivan_vysotin@bk.ru", "", "", "ivan_vysotin@bk.ru")
$ws.Range("C6").Style = $ws.Range("C2").Style

$null = $ws.Range("A7").Select()

# --- Page setup for printing ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
